# Apply the updates described by the diff:
#  - Metadata sheet: Version, Date, Contact fields
#  - Elements sheet: Binding Value Set for Device.code (row 14),
#    and Min / Base Min for Device.classCode (row 12)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z14").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAEntityCode"
$wsElem.Range("F12").Value = "0"
$wsElem.Range("AG12").Value = "0"
